$wb = $excel.ActiveWorkbook

# --- Features sheet: mark "New Things on the top" row as Fixed (blue fill) ---
$wsFeatures = $wb.Worksheets.Item("Features")
$wsFeatures.Range("B7:C7").Interior.Color = 15773696

# --- Issues sheet: append new issue rows ---
$wsIssues = $wb.Worksheets.Item("Issues")

$wsIssues.Cells.Item(25, 1).Value = 23

$wsIssues.Cells.Item(26, 1).Value = 24
$wsIssues.Cells.Item(26, 2).Value = "新的testrecord不会按降序添加"

$wsIssues.Cells.Item(27, 1).Value = 25
$wsIssues.Cells.Item(27, 5).Value = "SubPrograms全屏，ProgramEditView拉宽，app拉宽"
$wsIssues.Cells.Item(27, 2).Value = "UI调整"

$wsIssues.Cells.Item(28, 1).Value = 26

# --- Window/selection state: Features loses the active tab, Issues becomes active ---
$wsFeatures.Activate() | Out-Null
$wsFeatures.Range("F4").Select() | Out-Null

$wsIssues.Activate() | Out-Null
$wsIssues.Range("B33").Select() | Out-Null
